$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.775.79'
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").Value = '3.439.71'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''573.31'
$ws.Range("E5").Value = '  -1.12%  '

$ws.Range("D6").Value = '''159.02'
$ws.Range("E6").Value = '  -1.16%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.440.58'
$ws.Range("E8").Value = '  -1.47%  '

$ws.Range("D9").Value = '''0.581'
$ws.Range("E9").Value = '  -4.01%  '

$ws.Range("D10").Value = '''7.19'
$ws.Range("E10").Value = '  -1.58%  '

$ws.Range("E11").Value = '  -3.04%  '

$ws.Range("D12").Value = '''0.443'
$ws.Range("E12").Value = '  -1.17%  '

$ws.Range("D13").Value = '4.033.30'
$ws.Range("E13").Value = '  -1.46%  '

$ws.Range("E14").Value = '  -0.96%  '

$ws.Range("D15").Value = '''27.65'
$ws.Range("E15").Value = '  -3.70%  '

$ws.Range("D16").Value = '''0.0000180'
$ws.Range("E16").Value = '  -7.30%  '

$ws.Range("D17").Value = '64.795.60'
$ws.Range("E17").Value = '  -0.91%  '

$ws.Range("D18").Value = '3.440.83'
$ws.Range("E18").Value = '  -1.28%  '

$ws.Range("E19").Value = '  -2.26%  '

$ws.Range("D20").Value = '''13.89'
$ws.Range("E20").Value = '  -3.54%  '

$ws.Range("D21").Value = '''380.80'
$ws.Range("E21").Value = '  -1.65%  '

$ws.Range("D22").Value = '''7.96'
$ws.Range("E22").Value = '  -4.01%  '

$ws.Range("E23").Value = '  -1.31%  '

$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").Value = '''72.10'
$ws.Range("E25").Value = '  -1.61%  '

$ws.Range("E26").Value = '  -4.15%  '

$ws.Range("D27").Value = '''9.86'
$ws.Range("E27").Value = '  -1.64%  '

$ws.Range("E28").Value = '  -1.18%  '

$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("D31").Value = '''6.09'
$ws.Range("E31").Value = '  -2.92%  '

$ws.Range("E32").Value = '  -2.79%  '

$ws.Range("D33").Value = '''23.22'
$ws.Range("E33").Value = '  -1.93%  '

$ws.Range("E34").Value = '  -3.89%  '

$ws.Range("D35").Value = '''1.56'
$ws.Range("E35").Value = '  -1.25%  '

$ws.Range("D36").Value = '''161.40'
$ws.Range("E36").Value = '  -0.90%  '

$ws.Range("E37").Value = '  -2.13%  '

$ws.Range("D38").Value = '2.890.82'
$ws.Range("E38").Value = '  -3.66%  '

$ws.Range("D39").Value = '''0.0744'
$ws.Range("E39").Value = '  -4.56%  '

$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '''26.22'
$ws.Range("E40").Value = '  -3.81%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '''6.63'
$ws.Range("E41").Value = '  +0.88%  '

$ws.Range("D42").Value = '''4.53'
$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("D43").Value = '''42.91'
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").Value = '''0.0315'
$ws.Range("E44").Value = '  -2.95%  '

$ws.Range("E45").Value = '  -0.27%  '

$ws.Range("D46").Value = '''25.96'
$ws.Range("E46").Value = '  +1.22%  '

$ws.Range("E47").Value = '  +2.26%  '

$ws.Range("D48").Value = '''1.08'
$ws.Range("E48").Value = '  -3.03%  '

$ws.Range("D49").Value = '''315.69'
$ws.Range("E49").Value = '  -2.72%  '

$ws.Range("E50").Value = '  -3.75%  '

$ws.Range("E51").Value = '  -3.59%  '
